# "Hours and Task information updated !"
# Row 9 (Adit / Sr No 8) on Sheet1: the Interested Roles text gains a
# "Project Management " prefix and Hours per week goes from 5 to 6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D9").Value = "Project Management Requirements`nAna and design`nDev and V"
$ws.Range("E9").Value = 6

# Row heights re-flow around the edited/wrapped text.
$ws.Rows.Item(2).RowHeight = 48
$ws.Rows.Item(4).RowHeight = 32.25
$ws.Rows.Item(7).RowHeight = 32.25
$ws.Rows.Item(9).RowHeight = 63.75

# Scroll the view down a bit and leave the selection on the cell that was edited.
$excel.ActiveWindow.ScrollRow = 4
$ws.Range("E9").Select() | Out-Null
